$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.475.88"
$ws.Range("E2").Value = "  -2.40%  "

$ws.Range("D3").Value = "2.559.18"
$ws.Range("E3").Value = "  -3.81%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "520.21"
$ws.Range("E5").Value = "  -0.66%  "

$ws.Range("D6").Value = "143.60"
$ws.Range("E6").Value = "  -0.34%  "

$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("E8").Value = "  -1.46%  "

$ws.Range("D9").Value = "2.571.37"
$ws.Range("E9").Value = "  -3.62%  "

$ws.Range("D10").Value = "6.62"
$ws.Range("E10").Value = "  -4.30%  "

$ws.Range("D11").Value = "0.100"
$ws.Range("E11").Value = "  -2.37%  "

$ws.Range("D12").Value = "0.326"
$ws.Range("E12").Value = "  -2.95%  "

$ws.Range("E13").Value = "  -0.40%  "

$ws.Range("D14").Value = "3.014.17"
$ws.Range("E14").Value = "  -3.68%  "

$ws.Range("D15").Value = "57.446.44"
$ws.Range("E15").Value = "  -2.44%  "

$ws.Range("E16").Value = "  -4.32%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "0.0000133"
$ws.Range("E17").Value = "  -2.71%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.562.59"
$ws.Range("E18").Value = "  -4.23%  "

$ws.Range("D19").Value = "334.79"
$ws.Range("E19").Value = "  -1.28%  "

$ws.Range("D20").Value = "4.29"
$ws.Range("E20").Value = "  -2.33%  "

$ws.Range("D21").Value = "10.16"
$ws.Range("E21").Value = "  -2.39%  "

$ws.Range("D22").Value = "6.26"
$ws.Range("E22").Value = "  -1.41%  "

$ws.Range("D23").Value = "0.999"

$ws.Range("D24").Value = "64.69"
$ws.Range("E24").Value = "  +0.60%  "

$ws.Range("E25").Value = "  -0.74%  "

$ws.Range("B26").Value = "Polygon"
$ws.Range("C26").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D26").Value = "0.401"
$ws.Range("E26").Value = "  -4.94%  "

$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  -0.03%  "

$ws.Range("D28").Value = "2.679.33"
$ws.Range("E28").Value = "  -3.63%  "

$ws.Range("E29").Value = "  -2.98%  "

$ws.Range("D30").Value = "0.0₃0744"
$ws.Range("E30").Value = "  -7.55%  "

$ws.Range("E31").Value = "  +0.00%  "

$ws.Range("D32").Value = "6.22"
$ws.Range("E32").Value = "  -7.24%  "

$ws.Range("D34").Value = "18.59"
$ws.Range("E34").Value = "  -1.60%  "

$ws.Range("D35").Value = "148.63"
$ws.Range("E35").Value = "  -1.60%  "

$ws.Range("E36").Value = "  -3.25%  "

$ws.Range("E37").Value = "  -4.42%  "

$ws.Range("D38").Value = "0.839"
$ws.Range("E38").Value = "  -9.44%  "

$ws.Range("D39").Value = "35.99"
$ws.Range("E39").Value = "  -2.27%  "

$ws.Range("D40").Value = "0.830"
$ws.Range("E40").Value = "  -4.83%  "

$ws.Range("E41").Value = "  -2.13%  "

$ws.Range("E42").Value = "  -2.48%  "

$ws.Range("D44").Value = "267.79"
$ws.Range("E44").Value = "  -2.73%  "

$ws.Range("D45").Value = "10.66"
$ws.Range("E45").Value = "  -0.02%  "

$ws.Range("D46").Value = "0.0953"
$ws.Range("E46").Value = "  -1.52%  "

$ws.Range("E47").Value = "  -4.14%  "

$ws.Range("D48").Value = "18.78"
$ws.Range("E48").Value = "  -4.76%  "

$ws.Range("D49").Value = "0.0520"
$ws.Range("E49").Value = "  -2.61%  "

$ws.Range("D50").Value = "1.964.90"
$ws.Range("E50").Value = "  -4.72%  "

$ws.Range("D51").Value = "4.54"
$ws.Range("E51").Value = "  -4.03%  "
